$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the row-8 header labels before they get overwritten -------
# Before: A8=Date  E8=Acc.Code  F8=Cost Ctr.  G8=Account Title  I8=Debit  J8=Credit
# After : A8=Date  D8=Acc.Code  E8=Cost Ctr.  F8=Account Title  G8=Debit  H8=Credit
$accCode   = $ws.Range("E8").Value()
$costCtr   = $ws.Range("F8").Value()
$acctTitle = $ws.Range("G8").Value()
$debit     = $ws.Range("I8").Value()
$credit    = $ws.Range("J8").Value()

# --- Clear the cells that will no longer hold data (columns I, J) ------
$ws.Range("I7:J8").Clear()

# --- Write the header labels into their new (shifted-left) homes -------
$ws.Range("D8").Value = $accCode
$ws.Range("E8").Value = $costCtr
$ws.Range("F8").Value = $acctTitle
$ws.Range("G8").Value = $debit
$ws.Range("H8").Value = $credit

# --- Column widths -------------------------------------------------------
# New wide columns for the report content
$ws.Range("C1").EntireColumn.ColumnWidth = 61.166666666666664
$ws.Range("D1").EntireColumn.ColumnWidth = 17.5
$ws.Range("E1").EntireColumn.ColumnWidth = 17.166666666666668
$ws.Range("F1").EntireColumn.ColumnWidth = 54.833333333333336
# G (old col 7) is left untouched - keeps its existing bestFit width (12.5703125)
$ws.Range("H1").EntireColumn.ColumnWidth = 11.666666666666666

# Drop the now-unused trailing I/J column definitions in a single pass
$ws.Range("I1:J1").EntireColumn.Delete()

# --- Selection -------------------------------------------------------------
$ws.Range("C19").Select()
